# Apply the "maxGames"/"maxLateGames" numeric cleanup edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert E3 (maxLateGames for Minor) from text "6" to a real number 6.
$ws.Range("E3").Value = 6

# Column F holds "maxGames". Set every row (2-7) to the numeric value 16,
# then fix row 7 (ITB) back down to 8 - matching the new authoritative data.
$ws.Range("F2").Value = 16
$ws.Range("F3").Value = 16
$ws.Range("F4").Value = 16
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 16
$ws.Range("F7").Value = 8

# Move the active cell selection to J7 (matching the saved view state).
$ws.Range("J7").Select()
